# Commit: "Added geonet/knet readers, wrote sm2xml and ftpfetch programs."
#
# The underlying data (stations/locations/networks/coords) is unchanged.
# The only content edit is renaming the header cells for the coordinate
# columns from the verbose "latitude"/"longitude" to the shorter
# "lat"/"lon" (matching the new geonet/knet readers' column naming),
# and the active selection ends up on the new "lon" header cell (E1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "lat"
$ws.Range("E1").Value = "lon"

$ws.Range("E1").Select()
